$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 217 and 218 (the two "lámpara solar" products), shifting
# everything below up by two rows.
$ws.Rows("217:218").Delete()

# Excel recalculates the autofilter-backed hidden defined name
# (_xlnm._FilterDatabase) to track the table's new extent; make sure it
# reflects the now-smaller range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Hoja1!_FilterDatabase") {
        $n.RefersTo = "=Hoja1!`$A`$1:`$G`$233"
    }
}

# Leave the selection on the row that now occupies the first deleted
# row's position, matching the row-delete UI behavior.
$ws.Rows("217").Select()
